# Fixed update to excel issue
# 1. Rename "Requested quantity" headers to series-specific names.
# 2. Add a new "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$weeklySheet = $wb.Worksheets.Item("Weekly Quantity")
$weeklySheet.Range("B1").Value = "Weekly_PO_Qty"

$monthlySheet = $wb.Worksheets.Item("Monthly Trend")
$monthlySheet.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$poSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$poSheet.Name = "PO Forecast"

# Match header formatting (bold, centered, bordered) used on the other sheets.
$weeklySheet.Range("A1:B1").Copy()
$poSheet.Range("A1:D1").PasteSpecial(-4122)

# Match date-column formatting used on the other sheets.
$weeklySheet.Range("A2").Copy()
$poSheet.Range("A2:A48").PasteSpecial(-4122)

# Header row
$poSheet.Range("A1").Value = "ds"
$poSheet.Range("B1").Value = "PO_Forecast"
$poSheet.Range("C1").Value = "yhat_lower"
$poSheet.Range("D1").Value = "yhat_upper"

# Data rows
$poForecastData = @(
    @(45151.99999999999, 16, -12.8019972985016, 40.51178793944925),
    @(45186.99999999999, 16, -10.51339748316738, 43.01436179220958),
    @(45193.99999999999, 16, -9.026793995370337, 44.85927416262278),
    @(45207.99999999999, 17, -12.4489511328863, 41.89377118570111),
    @(45214.99999999999, 17, -9.407824150684984, 43.40862964692814),
    @(45221.99999999999, 17, -6.91145192835616, 43.9408981855724),
    @(45228.99999999999, 17, -9.725520258940648, 45.60024379568765),
    @(45235.99999999999, 17, -9.122997814533157, 43.70318303844179),
    @(45242.99999999999, 17, -9.231592475285602, 43.27783582815045),
    @(45270.99999999999, 17, -8.975565944087604, 45.79508058489585),
    @(45277.99999999999, 17, -9.150895773420068, 45.37621862838557),
    @(45298.99999999999, 18, -8.276864047119416, 45.68866739355133),
    @(45305.99999999999, 18, -9.989046147841666, 45.42589274809381),
    @(45312.99999999999, 18, -6.313657724685341, 43.83601231209287),
    @(45319.99999999999, 18, -9.826074732853158, 44.13179939035072),
    @(45326.99999999999, 18, -7.716730083134408, 42.25824972973549),
    @(45333.99999999999, 18, -7.234709198552871, 44.89567219160016),
    @(45340.99999999999, 18, -10.53927411687288, 43.95928719288867),
    @(45347.99999999999, 18, -6.347598816854497, 44.57288515878641),
    @(45368.99999999999, 19, -8.138931854762198, 42.94501686131439),
    @(45375.99999999999, 19, -8.892652100888508, 45.09850327696279),
    @(45382.99999999999, 19, -8.870304645070885, 44.38036200356646),
    @(45389.99999999999, 19, -10.14360010038757, 47.03730572421424),
    @(45459.99999999999, 20, -7.887166901103328, 47.2137499059675),
    @(45487.99999999999, 20, -5.611636317132308, 46.38547314684917),
    @(45501.99999999999, 20, -7.836082380333581, 47.50767221510329),
    @(45515.99999999999, 20, -6.414617958449852, 45.3928927420097),
    @(45529.99999999999, 21, -4.791590489282044, 46.51022063133515),
    @(45536.99999999999, 21, -7.513700671142912, 45.8950783058017),
    @(45543.99999999999, 21, -3.874596633828067, 47.96231175107339),
    @(45564.99999999999, 21, -5.904012100434492, 48.37455732044846),
    @(45571.99999999999, 21, -5.69401477161398, 46.01242098534638),
    @(45578.99999999999, 21, -4.321644340204176, 48.07639684146254),
    @(45585.99999999999, 21, -6.228396306028311, 49.15257887902686),
    @(45599.99999999999, 21, -5.927859900639453, 47.64905108807842),
    @(45606.99999999999, 22, -3.916277559764816, 47.80561054687085),
    @(45613.99999999999, 22, -5.089388364330618, 47.57689522108549),
    @(45634.99999999999, 22, -4.469199292301959, 48.63527370030954),
    @(45641.99999999999, 22, -4.835941012707272, 49.03932753446704),
    @(45648.99999999999, 22, -4.787526678838365, 52.14595282052802),
    @(45655.99999999999, 22, -3.311634807586605, 49.44355729508183),
    @(45662.99999999999, 22, -5.657587136832452, 47.79158919347555),
    @(45669.99999999999, 22, -5.231663382021282, 47.10467764435653),
    @(45676.99999999999, 22, -5.236396296141311, 47.87534306298378),
    @(45683.99999999999, 23, -2.947996327649104, 49.37468717924876),
    @(45690.99999999999, 23, -3.384259545474634, 50.24207768522551),
    @(45697.99999999999, 23, -3.116189779826009, 49.34897734233932)
)

for ($i = 0; $i -lt $poForecastData.Count; $i++) {
    $row = $poForecastData[$i]
    $r = $i + 2
    $poSheet.Cells.Item($r, 1).Value = $row[0]
    $poSheet.Cells.Item($r, 2).Value = $row[1]
    $poSheet.Cells.Item($r, 3).Value = $row[2]
    $poSheet.Cells.Item($r, 4).Value = $row[3]
}

$poSheet.Range("A1").Select()
